$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells (t="n") ---
$ws.Range("B7").Value = 35.10251782531287
$ws.Range("B19").Value = 23572221.57
$ws.Range("B20").Value = 8904762.379999999
$ws.Range("B24").Value = 720345
$ws.Range("B27").Value = 96046000
$ws.Range("B33").Value = 111282723.5587123
$ws.Range("B35").Value = 5245393.65
$ws.Range("B36").Value = 1125255.9
$ws.Range("B37").Value = 98468000
$ws.Range("B38").Value = 12160000
$ws.Range("B39").Value = 7025373.108712301
$ws.Range("B40").Value = 150000
$ws.Range("B41").Value = 5730000
$ws.Range("B43").Value = 5880000
$ws.Range("B47").Value = 4120000

# --- Text / percentage cells (t="inlineStr", stored as literal text) ---
# Force a Text number format first so the percent-looking strings are kept
# as literal text instead of being auto-converted to numeric percentages.
$textCells = @("B8", "B9", "B11", "B12", "B13", "B15", "B16", "B21", "B22", "B23", "B26")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B8").Value = "22.63%"
$ws.Range("B9").Value = "8.08%"
$ws.Range("B11").Value = "28.78%"
$ws.Range("B12").Value = "6.79%"
$ws.Range("B13").Value = "0.40%"
$ws.Range("B15").Value = "3.51%"
$ws.Range("B16").Value = "2.46%"
$ws.Range("B21").Value = "1.50%"
$ws.Range("B22").Value = "0.64%"
$ws.Range("B23").Value = "8.38%"
$ws.Range("B26").Value = "0.64%"
